# Add 2022-Q3 data:
#  - existing "2022-Q2" sheet is duplicated (the duplicate keeps the old Q2 numbers and
#    becomes the new, third "2022-Q2" sheet)
#  - the original sheet is renamed to "2022-Q3" and its figures are updated to the new quarter
#  - the "总计" (totals) sheet gets a new row for the 2022-Q3 entry

$wb = $excel.ActiveWorkbook

$q2 = $wb.Worksheets.Item("2022-Q2")

# 1. Duplicate the "2022-Q2" sheet right after itself; this copy keeps the old data
#    and will remain the "2022-Q2" sheet going forward.
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item($q2.Index + 1)

# 2. The original sheet becomes the 2022-Q3 sheet (rename it out of the way first so
#    the duplicate can reclaim the "2022-Q2" name).
$q3 = $q2
$q3.Name = "2022-Q3"
$q2Copy.Name = "2022-Q2"

# 3. Update the Q3 sheet's figures (fund code/name stay the same).
$q3.Range("Z1").Formula = "=""1.12"""
$q3.Range("AA1").Formula = "=""90.06"""
$q3.Range("AB1").Formula = "=""2.74"""
$q3.Range("AC1").Formula = "=""0.0307"""
$q3.Range("Z1:AC1").Copy()
$q3.Range("D2").PasteSpecial(-4163)
$q3.Range("Z1:AC1").ClearContents()
$q3.Range("H2").Value = 7

# Q3 sheet re-uses the "总计" sheet's header style for its header row / A2 cell.
$totals = $wb.Worksheets.Item("总计")
$totals.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$totals.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)

# 4. Update the totals sheet: existing 2022-Q2 row becomes 2022-Q3, and a new row is
#    appended for 2022-Q2 (mirroring the duplicated sheet order).
$totals.Range("B2").Value = "2022-Q3"

$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0.03
